# solved downloading issue of report id 172
# Update the various sample-data tables on the sheet (graph_excel.xlsx
# template) with the new figures used to regenerate the report charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shareholding pattern (rows 5-8) ---------------------------------
$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 70
$ws.Range("E5").Value = 70
$ws.Range("F5").Value = 70

$ws.Range("C6").Value = 5.85
$ws.Range("D6").Value = 11.93
$ws.Range("E6").Value = 10.8
$ws.Range("F6").Value = 11.18

$ws.Range("C7").Value = 0.01
$ws.Range("D7").Value = 0.16
$ws.Range("E7").Value = 0.16
$ws.Range("F7").Value = 0.15

$ws.Range("C8").Value = 24.22
$ws.Range("D8").Value = 15.61
$ws.Range("E8").Value = 16.95
$ws.Range("F8").Value = 18.86

# --- Liable to retire by rotation (rows 16 & 18) ----------------------
$ws.Range("C16").Value = 3
$ws.Range("C18").Value = 4

# --- ID / NID / SES style split (rows 28-29) --------------------------
$ws.Range("C28").Value = 0.5
$ws.Range("D28").Value = 0.5
$ws.Range("C29").Value = 0.5
$ws.Range("D29").Value = 0.5

# --- Year on Year Growth table (rows 38-42) ---------------------------
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 92.78

$ws.Range("C39").Value = 6.38
$ws.Range("D39").Value = 73.31

$ws.Range("C40").Value = 6.503
$ws.Range("D40").Value = 63.78

$ws.Range("C41").Value = 4.6016
$ws.Range("D41").Value = 83.6

$ws.Range("C42").Value = 5.2855
$ws.Range("D42").Value = 170.31

# --- Variation in Director's Remuneration (rows 50-51) ----------------
# Promoter row: the "NA" text is replaced with a numeric value.
$ws.Range("C50").Value = 4.24
$ws.Range("D50").Value = 0.16

# Non-Promoter row: label/value positions swap - the row now shows the
# "Non-Promoter" label followed by the "NA" text value.
$ws.Range("B51").Value = "Non-Promoter"
$ws.Range("C51").Value = "NA"
$ws.Range("D51").Value = 0.18

# --- Dividend / EPS / Payout table (rows 60-62) -----------------------
$ws.Range("C60").Value = 1.5
$ws.Range("D60").Value = 9
$ws.Range("E60").Value = 0.19

$ws.Range("C61").Value = 1.5
$ws.Range("D61").Value = 4.5
$ws.Range("E61").Value = 0.39

$ws.Range("C62").Value = 2
$ws.Range("D62").Value = 7.7
$ws.Range("E62").Value = 0.21

# --- Dividend / EPS / Payout (Today / 1Y / 3Y / 5Y) rows 72-74 --------
$ws.Range("B72").Value = 2
$ws.Range("C72").Value = 7.7
$ws.Range("D72").Value = 0.21

$ws.Range("B73").Value = 7
$ws.Range("C73").Value = 7.42
$ws.Range("D73").Value = 1.1

$ws.Range("B74").Value = 1.4
$ws.Range("C74").Value = 7.7
$ws.Range("D74").Value = 0.21

# --- Audit fees table (rows 83-85) ------------------------------------
$ws.Range("C83").Value = 0.38
$ws.Range("D83").Value = 0.25

$ws.Range("C84").Value = 0.3
$ws.Range("D84").Value = 0.3

$ws.Range("C85").Value = 0.02
$ws.Range("D85").Value = 0.05

# --- Audit fees table, second copy (rows 97-99) -----------------------
$ws.Range("C97").Value = 0.25
$ws.Range("D97").Value = 0.25
$ws.Range("E97").Value = 0.38

$ws.Range("C98").Value = 0.32
$ws.Range("D98").Value = 0.3
$ws.Range("E98").Value = 0.3

$ws.Range("C99").Value = 0.03
$ws.Range("D99").Value = 0.05
$ws.Range("E99").Value = 0.02
